# Applies the cryptos.xlsx price/volume/coin updates described in the commit diff.
# Values are written as plain text (matching the workbook's original inlineStr cells),
# using a temporary "@" (text) number format for purely-numeric-looking values so Excel
# does not silently convert them to numbers (which would drop significant trailing
# zeros, e.g. "0.0600" -> 0.06). The temporary format is reverted to the "Normal" style
# immediately afterwards so no visible/style differences remain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates (coin names, links, and already-textual price/volume strings) ---
$ws.Range("D2").Characters().Text = "35.651.08"
$ws.Range("E2").Characters().Text = "  -2.83%  "
$ws.Range("D3").Characters().Text = "1.985.39"
$ws.Range("E3").Characters().Text = "  -3.66%  "
$ws.Range("E4").Characters().Text = "  -0.05%  "
$ws.Range("E5").Characters().Text = "  +0.51%  "
$ws.Range("E6").Characters().Text = "  -4.64%  "
$ws.Range("E7").Characters().Text = "  +7.37%  "
$ws.Range("E8").Characters().Text = "  -0.01%  "
$ws.Range("E9").Characters().Text = "  -1.18%  "
$ws.Range("E10").Characters().Text = "  -0.91%  "
$ws.Range("E11").Characters().Text = "  -1.35%  "
$ws.Range("E12").Characters().Text = "  -2.25%  "
$ws.Range("E13").Characters().Text = "  +1.32%  "
$ws.Range("E14").Characters().Text = "  -0.84%  "
$ws.Range("D15").Characters().Text = "2.273.62"
$ws.Range("E15").Characters().Text = "  -3.69%  "
$ws.Range("E16").Characters().Text = "  -3.02%  "
$ws.Range("D17").Characters().Text = "1.975.44"
$ws.Range("E17").Characters().Text = "  -3.99%  "
$ws.Range("E18").Characters().Text = "  +7.90%  "
$ws.Range("D19").Characters().Text = "35.656.86"
$ws.Range("E19").Characters().Text = "  -2.59%  "
$ws.Range("E20").Characters().Text = "  -0.63%  "
$ws.Range("D21").Characters().Text = "0.0₃0852"
$ws.Range("E21").Characters().Text = "  -1.76%  "
$ws.Range("E22").Characters().Text = "  -1.17%  "
$ws.Range("E23").Characters().Text = "  -2.40%  "
$ws.Range("E24").Characters().Text = "  +0.05%  "
$ws.Range("E25").Characters().Text = "  +14.49%  "
$ws.Range("E26").Characters().Text = "  -4.59%  "
$ws.Range("E27").Characters().Text = "  -0.39%  "
$ws.Range("E28").Characters().Text = "  -0.16%  "
$ws.Range("E29").Characters().Text = "  -4.45%  "
$ws.Range("E30").Characters().Text = "  -2.67%  "
$ws.Range("E31").Characters().Text = "  -4.05%  "
$ws.Range("E32").Characters().Text = "  -5.57%  "
$ws.Range("E33").Characters().Text = "  +14.26%  "
$ws.Range("E34").Characters().Text = "  +0.16%  "
$ws.Range("E35").Characters().Text = "  +11.07%  "
$ws.Range("E36").Characters().Text = "  -2.71%  "
$ws.Range("E37").Characters().Text = "  -0.08%  "
$ws.Range("E38").Characters().Text = "  -2.64%  "
$ws.Range("E39").Characters().Text = "  +10.96%  "
$ws.Range("E40").Characters().Text = "  -1.53%  "
$ws.Range("E41").Characters().Text = "  -1.12%  "
$ws.Range("E42").Characters().Text = "  -2.76%  "
$ws.Range("E43").Characters().Text = "  +2.74%  "
$ws.Range("B44").Characters().Text = "ARBITRUM"
$ws.Range("C44").Characters().Text = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E44").Characters().Text = "  -1.23%  "
$ws.Range("B45").Characters().Text = "FraxShare"
$ws.Range("C45").Characters().Text = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("E45").Characters().Text = "  +3.36%  "
$ws.Range("B46").Characters().Text = "Aave"
$ws.Range("C46").Characters().Text = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E46").Characters().Text = "  -1.31%  "
$ws.Range("B47").Characters().Text = "InjectiveProtocol"
$ws.Range("C47").Characters().Text = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E47").Characters().Text = "  +1.48%  "
$ws.Range("D48").Characters().Text = "1.364.34"
$ws.Range("E48").Characters().Text = "  -3.20%  "
$ws.Range("E49").Characters().Text = "  -0.59%  "
$ws.Range("B50").Characters().Text = "RenderToken"
$ws.Range("C50").Characters().Text = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E50").Characters().Text = "  +1.05%  "
$ws.Range("B51").Characters().Text = "MultiversX"
$ws.Range("C51").Characters().Text = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("E51").Characters().Text = "  +3.68%  "

# --- Numeric-looking price updates: force text storage, then restore default styling ---
$numericRefs = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D14", "D18", "D20", "D22", "D23", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D50", "D51")
foreach ($ref in $numericRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D5").Value = "245.77"
$ws.Range("D6").Value = "0.636"
$ws.Range("D7").Value = "59.42"
$ws.Range("D9").Value = "59.34"
$ws.Range("D10").Value = "0.365"
$ws.Range("D11").Value = "0.0741"
$ws.Range("D12").Value = "0.105"
$ws.Range("D13").Value = "0.951"
$ws.Range("D14").Value = "14.73"
$ws.Range("D18").Value = "18.67"
$ws.Range("D20").Value = "71.73"
$ws.Range("D22").Value = "5.22"
$ws.Range("D23").Value = "233.41"
$ws.Range("D26").Value = "2.27"
$ws.Range("D27").Value = "9.24"
$ws.Range("D28").Value = "165.35"
$ws.Range("D29").Value = "19.30"
$ws.Range("D31").Value = "4.91"
$ws.Range("D32").Value = "1.13"
$ws.Range("D33").Value = "0.0962"
$ws.Range("D34").Value = "0.0600"
$ws.Range("D35").Value = "2.47"
$ws.Range("D36").Value = "4.38"
$ws.Range("D39").Value = "5.59"
$ws.Range("D40").Value = "1.24"
$ws.Range("D43").Value = "0.0930"
$ws.Range("D44").Value = "1.10"
$ws.Range("D45").Value = "7.79"
$ws.Range("D46").Value = "93.69"
$ws.Range("D47").Value = "16.39"
$ws.Range("D50").Value = "2.32"
$ws.Range("D51").Value = "47.27"

foreach ($ref in $numericRefs) {
    $ws.Range($ref).Style = "Normal"
}

